$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value2 = 132.11111
$ws.Range("J5").Value2 = 280
$ws.Range("L5").Value2 = 280
$ws.Range("N5").Value2 = -510

$ws.Range("H17").Value2 = 3520.5757
$ws.Range("I17").Value2 = 81
$ws.Range("J17").Value2 = 3742.484
$ws.Range("K17").Value2 = 243
$ws.Range("L17").Value2 = 11227.452
$ws.Range("M17").Value2 = -75
$ws.Range("N17").Value2 = -11563.452

$ws.Range("H51").Value2 = 3500.3333
$ws.Range("I51").Value2 = 2000.5
$ws.Range("J51").Value2 = 6500
$ws.Range("K51").Value2 = 2000.5
$ws.Range("L51").Value2 = 6500
$ws.Range("M51").Value2 = -1516.5
$ws.Range("N51").Value2 = -7468

$ws.Range("H62").Value2 = 4404.4443
$ws.Range("I62").Value2 = 2440
$ws.Range("J62").Value2 = 8333.333000000001
$ws.Range("K62").Value2 = 2440
$ws.Range("L62").Value2 = 8333.333000000001
$ws.Range("M62").Value2 = -1816
$ws.Range("N62").Value2 = -9581.333000000001

$ws.Range("H65").Value2 = 4404.4443
$ws.Range("I65").Value2 = 2440
$ws.Range("J65").Value2 = 8333.333000000001
$ws.Range("K65").Value2 = 12200
$ws.Range("L65").Value2 = 41666.665
$ws.Range("M65").Value2 = -9080
$ws.Range("N65").Value2 = -47906.665

$ws.Range("H94").Value2 = 2005
$ws.Range("I94").Value2 = 2005
$ws.Range("K94").Value2 = 2005
$ws.Range("M94").Value2 = -1554

$ws.Range("H98").Value2 = 2012.4783
$ws.Range("I98").Value2 = 1723.0588
$ws.Range("J98").Value2 = 2832.5
$ws.Range("K98").Value2 = 1723.0588
$ws.Range("L98").Value2 = 2832.5
$ws.Range("M98").Value2 = -225.0588
$ws.Range("N98").Value2 = -5828.5

$ws.Range("H100").Value2 = 1515.96
$ws.Range("I100").Value2 = 1363.5883
$ws.Range("J100").Value2 = 1839.75
$ws.Range("K100").Value2 = 1363.5883
$ws.Range("L100").Value2 = 1839.75
$ws.Range("M100").Value2 = -822.5882999999999
$ws.Range("N100").Value2 = -2921.75

$ws.Range("H106").Value2 = 2833.3333
$ws.Range("I106").Value2 = 2833.3333
$ws.Range("J106").Value2 = 0
$ws.Range("K106").Value2 = 2833.3333
$ws.Range("L106").Value2 = 0
$ws.Range("M106").Value2 = -2202.3333
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value2 = 482.1875
$ws.Range("I107").Value2 = 397.18518
$ws.Range("K107").Value2 = 397.18518
$ws.Range("M107").Value2 = 1522.81482

$ws.Range("H116").Value2 = 2800
$ws.Range("I116").Value2 = 2000
$ws.Range("K116").Value2 = 2000
$ws.Range("M116").Value2 = 1442

$ws.Range("H122").Value2 = 2012.4783
$ws.Range("I122").Value2 = 1723.0588
$ws.Range("J122").Value2 = 2832.5
$ws.Range("K122").Value2 = 5169.1764
$ws.Range("L122").Value2 = 8497.5
$ws.Range("M122").Value2 = -2719.1764
$ws.Range("N122").Value2 = -13397.5

$ws.Range("H135").Value2 = 93751820
$ws.Range("I135").Value2 = 45456924
$ws.Range("J135").Value2 = 200000620
$ws.Range("K135").Value2 = 409112316
$ws.Range("L135").Value2 = 1800005580
$ws.Range("M135").Value2 = -409109781
$ws.Range("N135").Value2 = -1800010650

$ws.Range("H137").Value2 = 1880.9025
$ws.Range("I137").Value2 = 1659.9062
$ws.Range("J137").Value2 = 2666.6667
$ws.Range("K137").Value2 = 4979.7186
$ws.Range("L137").Value2 = 8000.000100000001
$ws.Range("M137").Value2 = -2429.7186
$ws.Range("N137").Value2 = -13100.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 18344.27
$ws.Range("I32").Value2 = 20523.74
$ws.Range("K32").Value2 = 20523.74
$ws.Range("M32").Value2 = -20236.74

$ws.Range("H74").Value2 = 1879.25
$ws.Range("I74").Value2 = 1995.1154
$ws.Range("J74").Value2 = 1377.1666
$ws.Range("K74").Value2 = 1995.1154
$ws.Range("L74").Value2 = 1377.1666
$ws.Range("M74").Value2 = -1121.1154
$ws.Range("N74").Value2 = -3125.1666

$ws.Range("H77").Value2 = 1879.25
$ws.Range("I77").Value2 = 1995.1154
$ws.Range("J77").Value2 = 1377.1666
$ws.Range("K77").Value2 = 9975.576999999999
$ws.Range("L77").Value2 = 6885.833000000001
$ws.Range("M77").Value2 = -5607.576999999999
$ws.Range("N77").Value2 = -15621.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 1256.25
$ws.Range("I107").Value2 = 1007.6923
$ws.Range("J107").Value2 = 2333.3333
$ws.Range("K107").Value2 = 1007.6923
$ws.Range("L107").Value2 = 2333.3333
$ws.Range("M107").Value2 = 912.3077
$ws.Range("N107").Value2 = -6173.3333

$ws.Range("H135").Value2 = 0
$ws.Range("J135").Value2 = 0
$ws.Range("L135").Value2 = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 157.28572
$ws.Range("I7").Value2 = 105.77778
$ws.Range("K7").Value2 = 105.77778
$ws.Range("M7").Value2 = 7.222219999999993

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 206
$ws.Range("J23").Value2 = 206
$ws.Range("L23").Value2 = 618
$ws.Range("N23").Value2 = -1088

$ws.Range("H98").Value2 = 415.4
$ws.Range("I98").Value2 = 325
$ws.Range("J98").Value2 = 664
$ws.Range("K98").Value2 = 975
$ws.Range("L98").Value2 = 1992
$ws.Range("M98").Value2 = 523
$ws.Range("N98").Value2 = -4988

$ws.Range("H107").Value2 = 1161.862
$ws.Range("I107").Value2 = 340
$ws.Range("J107").Value2 = 1423.3636
$ws.Range("K107").Value2 = 1020
$ws.Range("L107").Value2 = 4270.0908
$ws.Range("M107").Value2 = 900
$ws.Range("N107").Value2 = -8110.0908

$ws.Range("H122").Value2 = 865.63635
$ws.Range("I122").Value2 = 586.3333
$ws.Range("J122").Value2 = 909.7368
$ws.Range("K122").Value2 = 5276.9997
$ws.Range("L122").Value2 = 8187.6312
$ws.Range("M122").Value2 = -2826.9997
$ws.Range("N122").Value2 = -13087.6312

$ws.Range("H131").Value2 = 1426.4546
$ws.Range("I131").Value2 = 1858
$ws.Range("J131").Value2 = 1066.8334
$ws.Range("K131").Value2 = 5574
$ws.Range("L131").Value2 = 3200.5002
$ws.Range("M131").Value2 = -534
$ws.Range("N131").Value2 = -13280.5002

$ws.Range("H134").Value2 = 3965.0967
$ws.Range("I134").Value2 = 3443
$ws.Range("K134").Value2 = 10329
$ws.Range("M134").Value2 = -5259

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 3357.7273
$ws.Range("I102").Value2 = 2923.2285
$ws.Range("J102").Value2 = 5047.4443
$ws.Range("K102").Value2 = 2923.2285
$ws.Range("L102").Value2 = 5047.4443
$ws.Range("M102").Value2 = -1301.2285
$ws.Range("N102").Value2 = -8291.444299999999

$ws.Range("H126").Value2 = 2679.12
$ws.Range("I126").Value2 = 1912.7142
$ws.Range("J126").Value2 = 3654.5454
$ws.Range("K126").Value2 = 5738.142599999999
$ws.Range("L126").Value2 = 10963.6362
$ws.Range("M126").Value2 = -3268.142599999999
$ws.Range("N126").Value2 = -15903.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 999.3333
$ws.Range("I96").Value2 = 799
$ws.Range("J96").Value2 = 1400
$ws.Range("K96").Value2 = 799
$ws.Range("L96").Value2 = 1400
$ws.Range("M96").Value2 = 574
$ws.Range("N96").Value2 = -4146

$ws.Range("H107").Value2 = 2323.4736
$ws.Range("I107").Value2 = 428
$ws.Range("J107").Value2 = 4429.5557
$ws.Range("K107").Value2 = 1284
$ws.Range("L107").Value2 = 13288.6671
$ws.Range("M107").Value2 = 636
$ws.Range("N107").Value2 = -17128.6671

$ws.Range("H113").Value2 = 601.5806
$ws.Range("I113").Value2 = 364.21054
$ws.Range("J113").Value2 = 977.4167
$ws.Range("K113").Value2 = 1092.63162
$ws.Range("L113").Value2 = 2932.2501
$ws.Range("M113").Value2 = 1077.36838
$ws.Range("N113").Value2 = -7272.2501
